$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "莎普爱思滴眼睛" (row 2) entirely; shifts subsequent rows up.
$ws.Rows.Item(2).Delete()
